$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold numeric-looking text values
# (e.g. "1.003", "0.1209", "  +4.79%  ") that must stay as text, not be
# auto-coerced into numbers by Excel smart entry. Pre-set the number
# format to Text ("@") for the whole data range before assigning values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.254.92'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.915.90'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '0.7439'
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").Value = '242.64'
$ws.Range("E6").Value = '  -3.09%  '
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = '0.3130'
$ws.Range("E8").Value = '  -2.76%  '
$ws.Range("D9").Value = '27.19'
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").Value = '0.06941'
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("D11").Value = '0.07993'
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").Value = '0.7683'
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").Value = '1.925.99'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = '5.284'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '91.23'
$ws.Range("E15").Value = '  -3.69%  '
$ws.Range("D16").Value = '30.305.64'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '14.17'
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").Value = '246.76'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("D19").Value = '5.823'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").Value = '0.000007840'
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = '2.186.59'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '6.593'
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("D25").Value = '9.370'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").Value = '165.09'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").Value = '18.89'
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("D28").Value = '0.1279'
$ws.Range("E28").Value = '  -3.57%  '
$ws.Range("D29").Value = '2.142'
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("D30").Value = '1.366'
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = '1.545'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").Value = '4.334'
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = '4.061'
$ws.Range("E33").Value = '  -2.30%  '
$ws.Range("D34").Value = '0.05172'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("D35").Value = '1.297'
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("D36").Value = '0.7446'
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '2.767'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '0.01933'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").Value = '2.770'
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").Value = '6.398'
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").Value = '75.68'
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("D42").Value = '0.4451'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '1.934'
$ws.Range("E43").Value = '  -2.83%  '
$ws.Range("D44").Value = '1.003'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '0.8362'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '101.26'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = '7.572'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '9.778'
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").Value = '37.05'
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1209'
$ws.Range("E50").Value = '  +4.79%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '952.73'
$ws.Range("E51").Value = '  -3.18%  '
